$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B (Coin name) updates for rows 44-46 (re-ranked coins) ---
$ws.Range("B44").Value = 'Mantle'
$ws.Range("B45").Value = 'OKB'
$ws.Range("B46").Value = 'ONDO'

# --- Column C (Link) updates for rows 44-46 (re-ranked coins) ---
$ws.Range("C44").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("C46").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'

# --- Column D (Price) updates ---
# Force text format so numeric-looking strings (e.g. "1.00", "69.40") are not
# auto-converted to numbers and lose their exact textual representation.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '60.119.61'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.340.18'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '564.91'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '130.50'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.340.59'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.912.58'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.343.33'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '24.66'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '60.236.01'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.49'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.14'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '355.19'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.475.75'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '69.40'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.67'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.48'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.371.44'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '22.97'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.40'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.90'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '158.81'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.40'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.750'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '40.79'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.19'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '23.65'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '22.39'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.893'

# --- Column E (Volume(1h) / % change) updates ---
$ws.Range("E2").Value = '  -5.60%  '
$ws.Range("E3").Value = '  -2.43%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("E5").Value = '  -2.45%  '
$ws.Range("E6").Value = '  +1.20%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("E8").Value = '  -2.43%  '
$ws.Range("E9").Value = '  -1.46%  '
$ws.Range("E10").Value = '  -1.88%  '
$ws.Range("E11").Value = '  -4.56%  '
$ws.Range("E12").Value = '  -1.10%  '
$ws.Range("E13").Value = '  -2.21%  '
$ws.Range("E14").Value = '  -0.22%  '
$ws.Range("E15").Value = '  -2.30%  '
$ws.Range("E16").Value = '  -3.69%  '
$ws.Range("E17").Value = '  -2.17%  '
$ws.Range("E18").Value = '  -5.42%  '
$ws.Range("E19").Value = '  +1.15%  '
$ws.Range("E20").Value = '  +1.33%  '
$ws.Range("E21").Value = '  -6.82%  '
$ws.Range("E22").Value = '  -7.15%  '
$ws.Range("E23").Value = '  -0.54%  '
$ws.Range("E24").Value = '  -2.34%  '
$ws.Range("E25").Value = '  -0.26%  '
$ws.Range("E26").Value = '  -6.22%  '
$ws.Range("E27").Value = '  +2.31%  '
$ws.Range("E28").Value = '  +19.09%  '
$ws.Range("E29").Value = '  +6.58%  '
$ws.Range("E30").Value = '  -0.07%  '
$ws.Range("E31").Value = '  +0.96%  '
$ws.Range("E32").Value = '  +0.88%  '
$ws.Range("E33").Value = '  -2.91%  '
$ws.Range("E34").Value = '  -0.06%  '
$ws.Range("E35").Value = '  -2.35%  '
$ws.Range("E36").Value = '  +1.22%  '
$ws.Range("E37").Value = '  +5.35%  '
$ws.Range("E38").Value = '  +2.83%  '
$ws.Range("E39").Value = '  +0.73%  '
$ws.Range("E40").Value = '  -3.20%  '
$ws.Range("E41").Value = '  +0.51%  '
$ws.Range("E42").Value = '  +0.09%  '
$ws.Range("E43").Value = '  +2.30%  '
$ws.Range("E44").Value = '  -4.26%  '
$ws.Range("E45").Value = '  -1.49%  '
$ws.Range("E46").Value = '  +8.05%  '
$ws.Range("E47").Value = '  +2.08%  '
$ws.Range("E48").Value = '  -1.18%  '
$ws.Range("E49").Value = '  +1.77%  '
$ws.Range("E50").Value = '  +10.59%  '
$ws.Range("E51").Value = '  +1.07%  '
